# The Devil's Faire - SlotMachineAssetDetails.xlsx
# "Added symbolTally obj literal, result tiles and power button"
#
# The RJust.* / LJust.* "result tile" font-size/justification labels used by
# rows 26-29 (CurrentBetText, WalletText, JackpotText, JackpotWinnerText) are
# being consolidated onto a shared "20px Arial" font label, and the
# right/left-justified coordinate labels are renamed/resized:
#   RJust.130x29 / RJust.721x29 / RJust.441x92 -> 20px Arial
#   LJust.42x29                                -> LJust.42x27
#   LJust.688x29                               -> LJust.688x27
#   RJust.345x92                               -> RJust.345x90

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 - WalletText result tile (touch this one first so the shared string
# table ends up in the same append order as the authoritative edit)
$ws.Range("E27").Value = "LJust.688x27"

# Row 26 - CurrentBetText result tile
$ws.Range("D26").Value = "20px Arial"
$ws.Range("E26").Value = "LJust.42x27"

$ws.Range("D27").Value = "20px Arial"

# Row 28 - JackpotText result tile
$ws.Range("D28").Value = "20px Arial"
$ws.Range("E28").Value = "RJust.345x90"

# Row 29 - JackpotWinnerText result tile
$ws.Range("D29").Value = "20px Arial"
$ws.Range("E29").Value = "RJust.345x90"

# Selection / view moved to E29, and the frozen/scrolled top-left cell
# reset back to A1 (topLeftCell attribute removed).
$ws.Range("E29").Select() | Out-Null
